$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / summary updates ---
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:48 PM"
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 25

# --- Insert 6 new rows for the Monday block (old row 32 "TOTAL" -> new row 38) ---
$ws.Range("A32:A37").EntireRow.Insert()

# --- Copy formatting (odd/even banding) into the 6 newly inserted blank rows ---
# Row 32/34 use the "odd" banding (style group from row 28); Row 33/35 use "even" banding (style group from row 29)
$ws.Range("A28:H28").Copy()
$ws.Range("A32:H32").PasteSpecial(-4122)
$ws.Range("A34:H34").PasteSpecial(-4122)
$ws.Range("A29:H29").Copy()
$ws.Range("A33:H33").PasteSpecial(-4122)
$ws.Range("A35:H35").PasteSpecial(-4122)
$ws.Range("A28:H28").Copy()
$ws.Range("A36:H36").PasteSpecial(-4122)
$ws.Range("A29:H29").Copy()
$ws.Range("A37:H37").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Write the full Monday + Tuesday line-item content (reordered, new items, all pricing zeroed) ---
$ws.Cells.Item(16, 1).Value = "Point 01"
$ws.Cells.Item(16, 2).Value = "CNA-TJ"
$ws.Cells.Item(16, 3).Value = "Inst"
$ws.Cells.Item(16, 4).Value = "CNA,Temporary Jumper"
$ws.Cells.Item(16, 5).Value = "EA"
$ws.Cells.Item(16, 6).Value = 3
$ws.Cells.Item(16, 8).Value = 0

$ws.Cells.Item(17, 1).Value = "Point 01"
$ws.Cells.Item(17, 2).Value = "SWI-15-CO1-100-H"
$ws.Cells.Item(17, 3).Value = "Inst"
$ws.Cells.Item(17, 4).Value = "SWI,15kV,Line Cutout 1PH,100A,Hook"
$ws.Cells.Item(17, 5).Value = "EA"
$ws.Cells.Item(17, 6).Value = 3
$ws.Cells.Item(17, 8).Value = 0

$ws.Cells.Item(18, 1).Value = "Point 02"
$ws.Cells.Item(18, 2).Value = "CNA-TJ"
$ws.Cells.Item(18, 3).Value = "Rem"
$ws.Cells.Item(18, 4).Value = "CNA,Temporary Jumper"
$ws.Cells.Item(18, 5).Value = "EA"
$ws.Cells.Item(18, 6).Value = 3
$ws.Cells.Item(18, 8).Value = 0

$ws.Cells.Item(19, 1).Value = "Point 02"
$ws.Cells.Item(19, 2).Value = "SWI-15-CO1-100-H"
$ws.Cells.Item(19, 3).Value = "Rem"
$ws.Cells.Item(19, 4).Value = "SWI,15kV,Line Cutout 1PH,100A,Hook"
$ws.Cells.Item(19, 5).Value = "EA"
$ws.Cells.Item(19, 6).Value = 3
$ws.Cells.Item(19, 8).Value = 0

$ws.Cells.Item(20, 1).Value = "Point 07"
$ws.Cells.Item(20, 2).Value = "ARM-8SF-GN-TL"
$ws.Cells.Item(20, 3).Value = "Inst"
$ws.Cells.Item(20, 4).Value = "ARM,8ft Sgl.Fiberglass,Gain,Tangent LD"
$ws.Cells.Item(20, 5).Value = "EA"
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 8).Value = 0

$ws.Cells.Item(21, 1).Value = "Point 07"
$ws.Cells.Item(21, 2).Value = "CNA-TR"
$ws.Cells.Item(21, 3).Value = "Inst"
$ws.Cells.Item(21, 4).Value = "CNA,Transfer Conductor"
$ws.Cells.Item(21, 5).Value = "EA"
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 8).Value = 0

$ws.Cells.Item(22, 1).Value = "Point 07"
$ws.Cells.Item(22, 2).Value = "CNA-TRPSE"
$ws.Cells.Item(22, 3).Value = "Inst"
$ws.Cells.Item(22, 4).Value = "CNA,Pole Set Energized Line,TempRelocate"
$ws.Cells.Item(22, 5).Value = "EA"
$ws.Cells.Item(22, 6).Value = 3
$ws.Cells.Item(22, 8).Value = 0

$ws.Cells.Item(23, 1).Value = "Point 07"
$ws.Cells.Item(23, 2).Value = "INS-15-P-S"
$ws.Cells.Item(23, 3).Value = "Inst"
$ws.Cells.Item(23, 4).Value = "INS,15kV,Pin,Silicon Polymer"
$ws.Cells.Item(23, 5).Value = "EA"
$ws.Cells.Item(23, 6).Value = 3
$ws.Cells.Item(23, 8).Value = 0

$ws.Cells.Item(24, 1).Value = "Point 07"
$ws.Cells.Item(24, 2).Value = "PIN-35-PTP"
$ws.Cells.Item(24, 3).Value = "Inst"
$ws.Cells.Item(24, 4).Value = "Pin,35kV,Pole Top"
$ws.Cells.Item(24, 5).Value = "EA"
$ws.Cells.Item(24, 6).Value = 1
$ws.Cells.Item(24, 8).Value = 0

$ws.Cells.Item(25, 1).Value = "Point 07"
$ws.Cells.Item(25, 2).Value = "PIN-XAL"
$ws.Cells.Item(25, 3).Value = "Inst"
$ws.Cells.Item(25, 4).Value = "Pin,Crossarm Light"
$ws.Cells.Item(25, 5).Value = "EA"
$ws.Cells.Item(25, 6).Value = 2
$ws.Cells.Item(25, 8).Value = 0

$ws.Cells.Item(26, 1).Value = "Point 07"
$ws.Cells.Item(26, 2).Value = "SAA-3-CV"
$ws.Cells.Item(26, 3).Value = "Inst"
$ws.Cells.Item(26, 4).Value = "SAA,3 inch,Clevis"
$ws.Cells.Item(26, 5).Value = "EA"
$ws.Cells.Item(26, 6).Value = 1
$ws.Cells.Item(26, 8).Value = 0

$ws.Cells.Item(27, 1).Value = "Point 07"
$ws.Cells.Item(27, 2).Value = "TIE-4-ALH-F"
$ws.Cells.Item(27, 3).Value = "Inst"
$ws.Cells.Item(27, 4).Value = "TIE,4 AWG,AL Hand Tie,F Neck"
$ws.Cells.Item(27, 5).Value = "EA"
$ws.Cells.Item(27, 6).Value = 24
$ws.Cells.Item(27, 8).Value = 0

$ws.Cells.Item(28, 1).Value = "Point 08"
$ws.Cells.Item(28, 2).Value = "ARM-8SF-GN-TL"
$ws.Cells.Item(28, 3).Value = "Rem"
$ws.Cells.Item(28, 4).Value = "ARM,8ft Sgl.Fiberglass,Gain,Tangent LD"
$ws.Cells.Item(28, 5).Value = "EA"
$ws.Cells.Item(28, 6).Value = 1
$ws.Cells.Item(28, 8).Value = 0

$ws.Cells.Item(29, 1).Value = "Point 08"
$ws.Cells.Item(29, 2).Value = "INS-15-P-S"
$ws.Cells.Item(29, 3).Value = "Rem"
$ws.Cells.Item(29, 4).Value = "INS,15kV,Pin,Silicon Polymer"
$ws.Cells.Item(29, 5).Value = "EA"
$ws.Cells.Item(29, 6).Value = 3
$ws.Cells.Item(29, 8).Value = 0

$ws.Cells.Item(30, 1).Value = "Point 08"
$ws.Cells.Item(30, 2).Value = "PIN-35-PTP"
$ws.Cells.Item(30, 3).Value = "Rem"
$ws.Cells.Item(30, 4).Value = "Pin,35kV,Pole Top"
$ws.Cells.Item(30, 5).Value = "EA"
$ws.Cells.Item(30, 6).Value = 1
$ws.Cells.Item(30, 8).Value = 0

$ws.Cells.Item(31, 1).Value = "Point 08"
$ws.Cells.Item(31, 2).Value = "PIN-XAL"
$ws.Cells.Item(31, 3).Value = "Rem"
$ws.Cells.Item(31, 4).Value = "Pin,Crossarm Light"
$ws.Cells.Item(31, 5).Value = "EA"
$ws.Cells.Item(31, 6).Value = 2
$ws.Cells.Item(31, 8).Value = 0

$ws.Cells.Item(32, 1).Value = "Point 08"
$ws.Cells.Item(32, 2).Value = "PLA-CUT"
$ws.Cells.Item(32, 3).Value = "Rem"
$ws.Cells.Item(32, 4).Value = "PLA,Cut Off Pole Top"
$ws.Cells.Item(32, 5).Value = "EA"
$ws.Cells.Item(32, 6).Value = 1
$ws.Cells.Item(32, 8).Value = 0

$ws.Cells.Item(33, 1).Value = "Point 08"
$ws.Cells.Item(33, 2).Value = "POL-40-2"
$ws.Cells.Item(33, 3).Value = "Rem"
$ws.Cells.Item(33, 4).Value = "Pole,40ft,Class 2"
$ws.Cells.Item(33, 5).Value = "EA"
$ws.Cells.Item(33, 6).Value = 1
$ws.Cells.Item(33, 8).Value = 0

$ws.Cells.Item(34, 1).Value = "Point 08"
$ws.Cells.Item(34, 2).Value = "SAA-3-CV"
$ws.Cells.Item(34, 3).Value = "Rem"
$ws.Cells.Item(34, 4).Value = "SAA,3 inch,Clevis"
$ws.Cells.Item(34, 5).Value = "EA"
$ws.Cells.Item(34, 6).Value = 1
$ws.Cells.Item(34, 8).Value = 0

$ws.Cells.Item(35, 1).Value = "Point 08"
$ws.Cells.Item(35, 2).Value = "TIE-4-ALH-F"
$ws.Cells.Item(35, 3).Value = "Rem"
$ws.Cells.Item(35, 4).Value = "TIE,4 AWG,AL Hand Tie,F Neck"
$ws.Cells.Item(35, 5).Value = "EA"
$ws.Cells.Item(35, 6).Value = 24
$ws.Cells.Item(35, 8).Value = 0

$ws.Cells.Item(36, 1).Value = "Point 03"
$ws.Cells.Item(36, 2).Value = "PLA-HDIG"
$ws.Cells.Item(36, 3).Value = "Inst"
$ws.Cells.Item(36, 4).Value = "PLA,Hand Dig or Additional  Excavation"
$ws.Cells.Item(36, 5).Value = "EA"
$ws.Cells.Item(36, 6).Value = 1
$ws.Cells.Item(36, 8).Value = 0

$ws.Cells.Item(37, 1).Value = "Point 05"
$ws.Cells.Item(37, 2).Value = "PLA-HDIG"
$ws.Cells.Item(37, 3).Value = "Inst"
$ws.Cells.Item(37, 4).Value = "PLA,Hand Dig or Additional  Excavation"
$ws.Cells.Item(37, 5).Value = "EA"
$ws.Cells.Item(37, 6).Value = 1
$ws.Cells.Item(37, 8).Value = 0

$ws.Cells.Item(43, 1).Value = "Point 05"
$ws.Cells.Item(43, 2).Value = "PLA-HDIG"
$ws.Cells.Item(43, 3).Value = "Inst"
$ws.Cells.Item(43, 4).Value = "PLA,Hand Dig or Additional  Excavation"
$ws.Cells.Item(43, 5).Value = "EA"
$ws.Cells.Item(43, 6).Value = 1
$ws.Cells.Item(43, 8).Value = 0

$ws.Cells.Item(44, 1).Value = "Point 07"
$ws.Cells.Item(44, 2).Value = "PLA-HDIG"
$ws.Cells.Item(44, 3).Value = "Inst"
$ws.Cells.Item(44, 4).Value = "PLA,Hand Dig or Additional  Excavation"
$ws.Cells.Item(44, 5).Value = "EA"
$ws.Cells.Item(44, 6).Value = 1
$ws.Cells.Item(44, 8).Value = 0

$ws.Cells.Item(45, 1).Value = "Point 09"
$ws.Cells.Item(45, 2).Value = "PLA-HDIG"
$ws.Cells.Item(45, 3).Value = "Inst"
$ws.Cells.Item(45, 4).Value = "PLA,Hand Dig or Additional  Excavation"
$ws.Cells.Item(45, 5).Value = "EA"
$ws.Cells.Item(45, 6).Value = 1
$ws.Cells.Item(45, 8).Value = 0

# --- Zero out the two TOTAL rows ---
$ws.Cells.Item(38, 8).Value = 0
$ws.Cells.Item(46, 8).Value = 0

